$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page 1")

# --- 1) Convert the "Temps" column (E) for the existing rows 3-26 from text
#     ("45 minutes") to plain numbers (45). The cell number format already
#     renders "0 minutes", so only the underlying value changes. We do NOT
#     touch the C/F text cells here so the shared-string table keeps its
#     original relative ordering (unused "X minutes" strings just get
#     garbage collected on save).
$times = @{
    3=45; 4=20; 5=100; 6=45; 7=10; 8=10; 9=30; 10=30; 11=30; 12=60;
    13=20; 14=20; 15=60; 16=20; 17=90; 18=30; 19=30; 20=60; 21=50;
    22=120; 23=120; 24=20; 25=60; 26=60
}
foreach ($row in ($times.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 5).Value = $times[$row]
}

# --- 2) Row 25 gains a "Description supplémentaire" (F) note.
$ws.Cells.Item(25, 6).Value = "Probleme de FK"

# --- 3) New journal entries, rows 27-35. Text cells are written in the
#     same order they were originally typed so the appended shared
#     strings land in the expected sequence.
$ws.Cells.Item(27, 3).Value = "Corriger un probleme de MLD"
$ws.Cells.Item(27, 4).Value = 44256
$ws.Cells.Item(27, 5).Value = 20
$ws.Cells.Item(27, 6).Value = "Tables intermediaire non relié + fk manquante"

$ws.Cells.Item(28, 3).Value = "Cas d'utilisation de la docs a remettre au propre"
$ws.Cells.Item(28, 4).Value = 44256
$ws.Cells.Item(28, 5).Value = 30
$ws.Cells.Item(28, 6).Value = "Comprehension de ce qu'on a discuter avec le professeur"

$ws.Cells.Item(29, 3).Value = "Creation de la page Creation de voyage"
$ws.Cells.Item(29, 4).Value = 44256
$ws.Cells.Item(29, 5).Value = 120

$ws.Cells.Item(30, 3).Value = "Premiere interaction avec la base de donnée"
$ws.Cells.Item(30, 4).Value = 44259
$ws.Cells.Item(30, 5).Value = 120

$ws.Cells.Item(31, 3).Value = "Codage du Login"
$ws.Cells.Item(31, 4).Value = 44259
$ws.Cells.Item(31, 5).Value = 120

$ws.Cells.Item(33, 6).Value = "Remplissage de la doc"
$ws.Cells.Item(33, 3).Value = "Finition du login"
$ws.Cells.Item(33, 4).Value = 44260
$ws.Cells.Item(33, 5).Value = 120

$ws.Cells.Item(32, 3).Value = "Bug du login a reparer"
$ws.Cells.Item(32, 4).Value = 44259
$ws.Cells.Item(32, 5).Value = 120

$ws.Cells.Item(30, 6).Value = "Gros bug sur la base de donnée et j'ai du changer de manière d'envoyer mes requetes SQL"
$ws.Cells.Item(32, 6).Value = "Suite au changement de l'approche de la gestion des requetes certains enorme bug sont aparue"

$ws.Cells.Item(34, 3).Value = "Profil"
$ws.Cells.Item(34, 4).Value = 44260
$ws.Cells.Item(34, 5).Value = 40

$ws.Cells.Item(35, 3).Value = "Remplisage journal de bord"
$ws.Cells.Item(35, 4).Value = 44260
$ws.Cells.Item(35, 5).Value = 20

# --- 4) Update the sheet's saved view state: scroll position and selection.
$ws.Range("F35").Select()
